$d = $word.ActiveDocument

# Locate the paragraph to remove by its text content, then expand the
# found range to the whole paragraph (including its paragraph mark) so
# that deleting it does not leave a stray empty/blank paragraph behind.
$rng = $d.Content
$found = $rng.Find.Execute(
    "There are two circuit design schemes we confirmed and it is hard to choose which one to use.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $rng.Expand(4) | Out-Null   # wdParagraph = 4
    $rng.Delete()
}
